# Update for release to deploy 0.1.1
# Remove the "ms-cmvgroup-observation" row (row 64) from the Observations
# sheet. Deleting the whole row shifts the subsequent rows (height,
# rh-status, weight) up by one and drops the now-unused shared strings.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(64).Delete()
